{"js": "// The CDC eviction moratorium cover letter was updated for the new 60-day\n// extension: the old date phrase \"July 31, 2021\" (built up out of several\n// runs) becomes \"October 3rd, 202\", and the \"_GoBack\" bookmark that used to\n// sit by itself in the blank paragraph right after this sentence moves up\n// to sit immediately after the new \"October 3rd\" text.\n\nconst body = context.document.body;\n\n// 1) Remove the existing \"_GoBack\" bookmark (currently alone in the blank\n//    paragraph right after the date sentence) so we can re-create it in its\n//    new location without leaving a duplicate behind.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Replace the old date text \"July 31, 2021\" with the new date text\n//    \"October 3rd, 202\", exactly as in the target revision. This collapses\n//    the several runs that used to spell out the old date into a single run.\nconst dateResults = body.search(\"July 31, 2021\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"October 3rd, 202\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Re-insert the \"_GoBack\" bookmark immediately after the new \"October\n//    3rd\" text (i.e. right before the \", 202\" that follows it).\nconst monthResults = body.search(\"October 3rd\", { matchCase: true });\nmonthResults.load(\"items\");\nawait context.sync();\n\nif (monthResults.items.length > 0) {\n  const afterMonth = monthResults.items[0].getRange(Word.RangeLocation.end);\n  afterMonth.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The CDC eviction moratorium cover letter was updated for the new 60-day\n# extension: the old date phrase \"July 31, 2021\" (built up out of several\n# runs) becomes \"October 3rd, 202\", and the \"_GoBack\" bookmark that used to\n# sit by itself in the blank paragraph right after this sentence moves up\n# to sit immediately after the new \"October 3rd\" text.\n\n$d = $word.ActiveDocument\n\n# 1) Replace the old date text \"July 31, 2021\" with the new date text\n#    \"October 3rd, 202\", exactly as in the target revision.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"July 31, 2021\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"October 3rd, 202\"\n$find.Execute(\n    $find.Text,             # FindText\n    $true,                  # MatchCase\n    $false,                 # MatchWholeWord\n    $false,                 # MatchWildcards\n    $false,                 # MatchSoundsLike\n    $false,                 # MatchAllWordForms\n    $true,                  # Forward\n    1,                      # Wrap (wdFindContinue)\n    $false,                 # Format\n    $find.Replacement.Text, # ReplaceWith\n    2                       # Replace (wdReplaceAll)\n)\n\n# 2) Re-insert the \"_GoBack\" bookmark immediately after the new \"October\n#    3rd\" text (right before the \", 202\" that follows it). Adding a\n#    bookmark with a name that already exists elsewhere in the document\n#    moves it here, so the old one (alone in the blank paragraph a couple\n#    of lines down) disappears from its old spot automatically.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"October 3rd\"\n$found2 = $find2.Execute()\nif ($found2) {\n    $afterMonth = $find2.Parent\n    $bookmarkPoint = $d.Range($afterMonth.End, $afterMonth.End)\n    $d.Bookmarks.Add(\"_GoBack\", $bookmarkPoint)\n}\n"}
